# Remove the auto-advancing timer (and the slow-speed slide transition that
# rode along with it) from every slide that has one. In the OOXML this is the
# <mc:AlternateContent> wrapper around <p:transition spd="slow" ... advTm="..."/>
# at the end of <p:sld> — clearing the transition's Duration drops the whole
# <p:transition>/<mc:AlternateContent> element instead of merely zeroing an
# attribute.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $t = $s.SlideShowTransition

    if ($t.AdvanceOnTime) {
        $t.Duration = $null
    }
}
